$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 239; $i++) {
    $ws.Cells.Item($i, 3).Value = 45203
}
